$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.000005906118310345425
$ws.Range("E2").Value = 0.000005906118310345425

# Row 3
$ws.Range("D3").Value = 0.9999999945753031
$ws.Range("E3").Value = 0.9999999945753031

# Row 4
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = 0.1085728188178906
$ws.Range("E4").Value = 0.1085728188178906

# Row 5
$ws.Range("D5").Value = 0.0000000000000000000000000000000000003785093707568907
$ws.Range("E5").Value = 0.0000000000000000000000000000000000003785093707568907

# Row 6
$ws.Range("D6").Value = 0.002049544324786162
$ws.Range("E6").Value = 0.002049544324786162

# Row 7
$ws.Range("D7").Value = 0.9999999992763677
$ws.Range("E7").Value = 0.0000000007236322652204308

# Row 8
$ws.Range("D8").Value = 0.9990833781283756
$ws.Range("E8").Value = 0.0009166218716244456

# Row 9
$ws.Range("D9").Value = 0.9999999999998532
$ws.Range("E9").Value = 0.0000000000001467714838554457

# Row 10
$ws.Range("D10").Value = 0.9998675295770668
$ws.Range("E10").Value = 0.0001324704229331841

# Row 11
$ws.Range("D11").Value = 0.9999999987548187
$ws.Range("E11").Value = 0.000000001245181291409381
$ws.Range("F11").Value = 1.915034532546997
$ws.Range("G11").Value = 0.9
